# PCB BOM redesign: audio connector changed from the Lumberg/Reichelt part
# (P104 / LUM150302 / "Audio Out") to the Thonk PJ301M-12 jack (J101 /
# Thonkicon / "AUOUT"), plus a new datasheet/shop-link note in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 holds the audio-out connector BOM line - swap its contents to the
# new part.
$ws.Range("A10").Value = "Thonk"
$ws.Range("C10").Value = "J101"
$ws.Range("D10").Value = "Thonkicon"
$ws.Range("F10").Value = "AUOUT"
$ws.Range("G10").Value = "PJ301M-12"
$ws.Range("I10").Value = "https://www.thonk.co.uk/shop/3-5mm-jacks/"

# Leave the selection on the newly added note cell, matching where the
# author ended up after editing.
$ws.Range("I10").Select()
